$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“핸즈온 LLM” 번역 완료"
$ws.Range("E12").Value = "https://tensorflow.blog/2025/03/29/%ed%95%b8%ec%a6%88%ec%98%a8-llm-%eb%b2%88%ec%97%ad-%ec%99%84%eb%a3%8c/"
